# Auto-generated edit script applying the Tiamat_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 2634911.8
$ws.Range("I74").Value = 3128315
$ws.Range("J74").Value = 3428.3333
$ws.Range("K74").Value = 3128315
$ws.Range("L74").Value = 3428.3333
$ws.Range("M74").Value = -3127379
$ws.Range("N74").Value = -5300.3333
$ws.Range("H76").Value = 14928546
$ws.Range("I76").Value = 16396357
$ws.Range("J76").Value = 5794.6665
$ws.Range("K76").Value = 16396357
$ws.Range("L76").Value = 5794.6665
$ws.Range("M76").Value = -16396042
$ws.Range("N76").Value = -6424.6665
$ws.Range("H77").Value = 2634911.8
$ws.Range("I77").Value = 3128315
$ws.Range("J77").Value = 3428.3333
$ws.Range("K77").Value = 15641575
$ws.Range("L77").Value = 17141.6665
$ws.Range("M77").Value = -15636895
$ws.Range("N77").Value = -26501.6665
$ws.Range("H79").Value = 14928546
$ws.Range("I79").Value = 16396357
$ws.Range("J79").Value = 5794.6665
$ws.Range("K79").Value = 16396357
$ws.Range("L79").Value = 5794.6665
$ws.Range("M79").Value = -16395265
$ws.Range("N79").Value = -7978.6665
$ws.Range("H80").Value = 5253605
$ws.Range("I80").Value = 6667675.5
$ws.Range("J80").Value = 4546569.5
$ws.Range("K80").Value = 20003026.5
$ws.Range("L80").Value = 13639708.5
$ws.Range("M80").Value = -20002028.5
$ws.Range("N80").Value = -13641704.5
$ws.Range("H83").Value = 5253605
$ws.Range("I83").Value = 6667675.5
$ws.Range("J83").Value = 4546569.5
$ws.Range("K83").Value = 60009079.5
$ws.Range("L83").Value = 40919125.5
$ws.Range("M83").Value = -60004087.5
$ws.Range("N83").Value = -40929109.5
$ws.Range("H92").Value = 20834466
$ws.Range("I92").Value = 29412918
$ws.Range("J92").Value = 1085.7142
$ws.Range("K92").Value = 29412918
$ws.Range("L92").Value = 1085.7142
$ws.Range("M92").Value = -29411670
$ws.Range("N92").Value = -3581.7142
$ws.Range("H100").Value = 3712.2
$ws.Range("I100").Value = 2444.9443
$ws.Range("J100").Value = 6970.857
$ws.Range("K100").Value = 2444.9443
$ws.Range("L100").Value = 6970.857
$ws.Range("M100").Value = -1903.9443
$ws.Range("N100").Value = -8052.857
$ws.Range("H103").Value = 1602814.2
$ws.Range("I103").Value = 1831756.2
$ws.Range("J103").Value = 220
$ws.Range("K103").Value = 5495268.6
$ws.Range("L103").Value = 660
$ws.Range("M103").Value = -5494682.6
$ws.Range("N103").Value = -1832
$ws.Range("H131").Value = 5312.143
$ws.Range("I131").Value = 1046.25
$ws.Range("J131").Value = 11000
$ws.Range("K131").Value = 3138.75
$ws.Range("L131").Value = 33000
$ws.Range("M131").Value = 1901.25
$ws.Range("N131").Value = -43080
$ws.Range("H138").Value = 1850.47
$ws.Range("I138").Value = 710.08
$ws.Range("J138").Value = 2230.6
$ws.Range("K138").Value = 2130.24
$ws.Range("L138").Value = 6691.799999999999
$ws.Range("M138").Value = 3009.76
$ws.Range("N138").Value = -16971.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 209480.86
$ws.Range("I32").Value = 260388.2
$ws.Range("J32").Value = 72667.375
$ws.Range("K32").Value = 260388.2
$ws.Range("L32").Value = 72667.375
$ws.Range("M32").Value = -260101.2
$ws.Range("N32").Value = -73241.375
$ws.Range("H97").Value = 1458.2059
$ws.Range("I97").Value = 710.6539
$ws.Range("J97").Value = 3887.75
$ws.Range("K97").Value = 710.6539
$ws.Range("L97").Value = 3887.75
$ws.Range("M97").Value = -214.6539
$ws.Range("N97").Value = -4879.75
$ws.Range("H122").Value = 949
$ws.Range("I122").Value = 832.86957
$ws.Range("J122").Value = 1330.5714
$ws.Range("K122").Value = 2498.60871
$ws.Range("L122").Value = 3991.7142
$ws.Range("M122").Value = -48.60870999999997
$ws.Range("N122").Value = -8891.7142
$ws.Range("H132").Value = 9479191
$ws.Range("I132").Value = 15012581
$ws.Range("J132").Value = 564284.75
$ws.Range("K132").Value = 45037743
$ws.Range("L132").Value = 1692854.25
$ws.Range("M132").Value = -45035213
$ws.Range("N132").Value = -1697914.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 10326.143
$ws.Range("I94").Value = 11141.454
$ws.Range("K94").Value = 11141.454
$ws.Range("M94").Value = -10690.454
$ws.Range("H134").Value = 14316316
$ws.Range("I134").Value = 26317296
$ws.Range("K134").Value = 78951888
$ws.Range("M134").Value = -78949353

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1035.6
$ws.Range("I16").Value = 989.5
$ws.Range("J16").Value = 1066.3334
$ws.Range("K16").Value = 989.5
$ws.Range("L16").Value = 1066.3334
$ws.Range("M16").Value = -702.5
$ws.Range("N16").Value = -1640.3334
$ws.Range("H31").Value = 5904.041
$ws.Range("I31").Value = 6073.7
$ws.Range("J31").Value = 5787.0347
$ws.Range("K31").Value = 6073.7
$ws.Range("L31").Value = 5787.0347
$ws.Range("M31").Value = -5778.7
$ws.Range("N31").Value = -6377.0347
$ws.Range("H34").Value = 5904.041
$ws.Range("I34").Value = 6073.7
$ws.Range("J34").Value = 5787.0347
$ws.Range("K34").Value = 6073.7
$ws.Range("L34").Value = 5787.0347
$ws.Range("M34").Value = -5871.7
$ws.Range("N34").Value = -6191.0347
$ws.Range("H113").Value = 1035.6
$ws.Range("I113").Value = 989.5
$ws.Range("J113").Value = 1066.3334
$ws.Range("K113").Value = 989.5
$ws.Range("L113").Value = 1066.3334
$ws.Range("M113").Value = 1180.5
$ws.Range("N113").Value = -5406.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 689.5
$ws.Range("I117").Value = 200
$ws.Range("J117").Value = 715.2632
$ws.Range("K117").Value = 600
$ws.Range("L117").Value = 2145.7896
$ws.Range("M117").Value = 2842
$ws.Range("N117").Value = -9029.7896

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 45765
$ws.Range("J49").Value = 45765
$ws.Range("L49").Value = 45765
$ws.Range("N49").Value = -46133
$ws.Range("H70").Value = 4047.8696
$ws.Range("I70").Value = 3933.5557
$ws.Range("J70").Value = 4121.357
$ws.Range("K70").Value = 3933.5557
$ws.Range("L70").Value = 4121.357
$ws.Range("M70").Value = -3663.5557
$ws.Range("N70").Value = -4661.357
$ws.Range("H73").Value = 4047.8696
$ws.Range("I73").Value = 3933.5557
$ws.Range("J73").Value = 4121.357
$ws.Range("K73").Value = 3933.5557
$ws.Range("L73").Value = 4121.357
$ws.Range("M73").Value = -2997.5557
$ws.Range("N73").Value = -5993.357
$ws.Range("H97").Value = 1017.6923
$ws.Range("I97").Value = 929
$ws.Range("J97").Value = 1505.5
$ws.Range("K97").Value = 929
$ws.Range("L97").Value = 1505.5
$ws.Range("M97").Value = -433
$ws.Range("N97").Value = -2497.5
$ws.Range("H107").Value = 399.45
$ws.Range("I107").Value = 346.41177
$ws.Range("J107").Value = 700
$ws.Range("K107").Value = 346.41177
$ws.Range("L107").Value = 700
$ws.Range("M107").Value = 1573.58823
$ws.Range("N107").Value = -4540
$ws.Range("H132").Value = 41964.27
$ws.Range("I132").Value = 2584.7
$ws.Range("K132").Value = 7754.099999999999
$ws.Range("M132").Value = -5224.099999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2571.5557
$ws.Range("I7").Value = 1806
$ws.Range("J7").Value = 3184
$ws.Range("K7").Value = 1806
$ws.Range("L7").Value = 3184
$ws.Range("M7").Value = -1694
$ws.Range("N7").Value = -3408
$ws.Range("H42").Value = 43997.777
$ws.Range("J42").Value = 43997.777
$ws.Range("L42").Value = 43997.777
$ws.Range("N42").Value = -45123.777
$ws.Range("H46").Value = 251
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 251
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 251
$ws.Range("M46").Value = $null
$ws.Range("N46").Value = -627
$ws.Range("H49").Value = 43997.777
$ws.Range("J49").Value = 43997.777
$ws.Range("L49").Value = 43997.777
$ws.Range("N49").Value = -44291.777
$ws.Range("H50").Value = 16200
$ws.Range("J50").Value = 16200
$ws.Range("L50").Value = 16200
$ws.Range("N50").Value = -17474
$ws.Range("H54").Value = 33400
$ws.Range("J54").Value = 33400
$ws.Range("L54").Value = 33400
$ws.Range("N54").Value = -34688
$ws.Range("H126").Value = 2571.5557
$ws.Range("I126").Value = 1806
$ws.Range("J126").Value = 3184
$ws.Range("K126").Value = 5418
$ws.Range("L126").Value = 9552
$ws.Range("M126").Value = -2948
$ws.Range("N126").Value = -14492

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3090.0264
$ws.Range("I122").Value = 2261.8696
$ws.Range("J122").Value = 4359.8667
$ws.Range("K122").Value = 6785.6088
$ws.Range("L122").Value = 13079.6001
$ws.Range("M122").Value = -4335.6088
$ws.Range("N122").Value = -17979.6001
